# Auto-generated edit script applying the Siren_Profits.xlsx market-data refresh diff.
# Updates literal numeric cells (currentAveragePrice / Leve price & profit columns)
# across all 8 Leve sheets, matching the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 5496.4287
$ws.Range("J48").Value = 5496.4287
$ws.Range("L48").Value = 16489.2861
$ws.Range("N48").Value = -17073.2861

$ws.Range("H56").Value = 5496.4287
$ws.Range("J56").Value = 5496.4287
$ws.Range("L56").Value = 16489.2861
$ws.Range("N56").Value = -17557.2861

$ws.Range("H58").Value = 2985.818
$ws.Range("J58").Value = 4083.1667
$ws.Range("L58").Value = 12249.5001
$ws.Range("N58").Value = -12549.5001

$ws.Range("H137").Value = 731099.3
$ws.Range("I137").Value = 982837.5
$ws.Range("K137").Value = 2948512.5
$ws.Range("M137").Value = -2945962.5

$ws.Range("H138").Value = 223863.94
$ws.Range("I138").Value = 600215.4
$ws.Range("J138").Value = 5976.263
$ws.Range("K138").Value = 1800646.2
$ws.Range("L138").Value = 17928.789
$ws.Range("M138").Value = -1795506.2
$ws.Range("N138").Value = -28208.789


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3943.15
$ws.Range("I32").Value = 3943.15
$ws.Range("K32").Value = 3943.15
$ws.Range("M32").Value = -3656.15

$ws.Range("H45").Value = 83076.80499999999
$ws.Range("I45").Value = 116894.28
$ws.Range("K45").Value = 116894.28
$ws.Range("M45").Value = -116517.28

$ws.Range("H61").Value = 5074.7188
$ws.Range("I61").Value = 5654.864
$ws.Range("J61").Value = 3798.4
$ws.Range("K61").Value = 5654.864
$ws.Range("L61").Value = 3798.4
$ws.Range("M61").Value = -5442.864
$ws.Range("N61").Value = -4222.4

$ws.Range("H74").Value = 3747.3389
$ws.Range("I74").Value = 15911.5
$ws.Range("K74").Value = 15911.5
$ws.Range("M74").Value = -15037.5

$ws.Range("H77").Value = 3747.3389
$ws.Range("I77").Value = 15911.5
$ws.Range("K77").Value = 79557.5
$ws.Range("M77").Value = -75189.5

$ws.Range("H122").Value = 973615.8
$ws.Range("I122").Value = 5891.4585
$ws.Range("K122").Value = 17674.3755
$ws.Range("M122").Value = -15224.3755

$ws.Range("H132").Value = 2163.2693
$ws.Range("I132").Value = 1793.2609
$ws.Range("K132").Value = 5379.7827
$ws.Range("M132").Value = -2849.7827

$ws.Range("H136").Value = 5074.7188
$ws.Range("I136").Value = 5654.864
$ws.Range("J136").Value = 3798.4
$ws.Range("K136").Value = 16964.592
$ws.Range("L136").Value = 11395.2
$ws.Range("M136").Value = -14414.592
$ws.Range("N136").Value = -16495.2


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8652.846
$ws.Range("I86").Value = 13516.667
$ws.Range("J86").Value = 4483.857
$ws.Range("K86").Value = 13516.667
$ws.Range("L86").Value = 4483.857
$ws.Range("M86").Value = -12393.667
$ws.Range("N86").Value = -6729.857

$ws.Range("H89").Value = 8652.846
$ws.Range("I89").Value = 13516.667
$ws.Range("J89").Value = 4483.857
$ws.Range("K89").Value = 67583.33499999999
$ws.Range("L89").Value = 22419.285
$ws.Range("M89").Value = -61967.33499999999
$ws.Range("N89").Value = -33651.285

$ws.Range("H107").Value = 1438.75
$ws.Range("I107").Value = 1501.4286
$ws.Range("K107").Value = 1501.4286
$ws.Range("M107").Value = 418.5714

$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 4000
$ws.Range("K17").Value = 4000
$ws.Range("M17").Value = -3826

$ws.Range("H31").Value = 4791
$ws.Range("I31").Value = 1522
$ws.Range("J31").Value = 6098.6
$ws.Range("K31").Value = 1522
$ws.Range("L31").Value = 6098.6
$ws.Range("M31").Value = -1227
$ws.Range("N31").Value = -6688.6

$ws.Range("H34").Value = 4791
$ws.Range("I34").Value = 1522
$ws.Range("J34").Value = 6098.6
$ws.Range("K34").Value = 1522
$ws.Range("L34").Value = 6098.6
$ws.Range("M34").Value = -1320
$ws.Range("N34").Value = -6502.6

$ws.Range("H38").Value = 3115.75
$ws.Range("I38").Value = 3012.6667
$ws.Range("J38").Value = 3425
$ws.Range("K38").Value = 3012.6667
$ws.Range("L38").Value = 3425
$ws.Range("M38").Value = -2635.6667
$ws.Range("N38").Value = -4179

$ws.Range("H46").Value = 3115.75
$ws.Range("I46").Value = 3012.6667
$ws.Range("J46").Value = 3425
$ws.Range("K46").Value = 3012.6667
$ws.Range("L46").Value = 3425
$ws.Range("M46").Value = -2801.6667
$ws.Range("N46").Value = -3847

$ws.Range("H132").Value = 7077.0527
$ws.Range("I132").Value = 8031
$ws.Range("K132").Value = 24093
$ws.Range("M132").Value = -21563

$ws.Range("H141").Value = 411738.88
$ws.Range("J141").Value = 517930.94
$ws.Range("L141").Value = 517930.94
$ws.Range("N141").Value = -528290.9399999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2058.7273
$ws.Range("I25").Value = 118
$ws.Range("J25").Value = 2490
$ws.Range("K25").Value = 354
$ws.Range("L25").Value = 7470
$ws.Range("M25").Value = -185
$ws.Range("N25").Value = -7808

$ws.Range("H30").Value = 2058.7273
$ws.Range("I30").Value = 118
$ws.Range("J30").Value = 2490
$ws.Range("K30").Value = 354
$ws.Range("L30").Value = 7470
$ws.Range("M30").Value = -252
$ws.Range("N30").Value = -7674

$ws.Range("H58").Value = 2947.5264
$ws.Range("I58").Value = 503
$ws.Range("J58").Value = 3083.3333
$ws.Range("K58").Value = 1509
$ws.Range("L58").Value = 9249.999899999999
$ws.Range("M58").Value = -1381
$ws.Range("N58").Value = -9505.999899999999

$ws.Range("H121").Value = 1807.3077
$ws.Range("I121").Value = 965.6667
$ws.Range("J121").Value = 2528.7144
$ws.Range("K121").Value = 2897.0001
$ws.Range("L121").Value = 7586.1432
$ws.Range("M121").Value = -1587.0001
$ws.Range("N121").Value = -10206.1432

$ws.Range("H122").Value = 2856.6785
$ws.Range("J122").Value = 3604.4443
$ws.Range("L122").Value = 32439.9987
$ws.Range("N122").Value = -37339.9987

$ws.Range("H137").Value = 8425.177
$ws.Range("J137").Value = 9863.77
$ws.Range("L137").Value = 29591.31
$ws.Range("N137").Value = -39791.31


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 782.5
$ws.Range("I107").Value = 844.2222
$ws.Range("J107").Value = 597.3333
$ws.Range("K107").Value = 844.2222
$ws.Range("L107").Value = 597.3333
$ws.Range("M107").Value = 1075.7778
$ws.Range("N107").Value = -4437.3333

$ws.Range("H113").Value = 21163.5
$ws.Range("I113").Value = 29245.25
$ws.Range("K113").Value = 29245.25
$ws.Range("M113").Value = -27075.25

$ws.Range("H126").Value = 34148.8
$ws.Range("J126").Value = 15915
$ws.Range("L126").Value = 47745
$ws.Range("N126").Value = -52685


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1419.8
$ws.Range("I46").Value = 846.0769
$ws.Range("J46").Value = 2485.2856
$ws.Range("K46").Value = 846.0769
$ws.Range("L46").Value = 2485.2856
$ws.Range("M46").Value = -658.0769
$ws.Range("N46").Value = -2861.2856

$ws.Range("H122").Value = 6766.1924
$ws.Range("I122").Value = 6744.4287
$ws.Range("K122").Value = 20233.2861
$ws.Range("M122").Value = -17783.2861

$ws.Range("H136").Value = 8628.375
$ws.Range("I136").Value = 6041.357
$ws.Range("K136").Value = 18124.071
$ws.Range("M136").Value = -15574.071


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 88499.5
$ws.Range("J118").Value = 88499.5
$ws.Range("L118").Value = 88499.5
$ws.Range("N118").Value = -91813.5

$ws.Range("H129").Value = 95000
$ws.Range("J129").Value = 95000
$ws.Range("L129").Value = 95000
$ws.Range("N129").Value = -105000

$ws.Range("H132").Value = 8091.811
$ws.Range("I132").Value = 9874.673000000001
$ws.Range("K132").Value = 29624.019
$ws.Range("M132").Value = -27094.019

